$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.693.48"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "1.888.30"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "235.74"
$ws.Range("E5").Value = "  -4.03%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "0.4874"
$ws.Range("E7").Value = "  -2.17%  "

$ws.Range("D8").Value = "0.2893"
$ws.Range("E8").Value = "  -3.36%  "

$ws.Range("D9").Value = "0.06660"
$ws.Range("E9").Value = "  -3.17%  "

$ws.Range("D10").Value = "1.883.96"
$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("D11").Value = "16.62"
$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("D12").Value = "0.07230"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").Value = "89.10"
$ws.Range("E13").Value = "  -2.09%  "

$ws.Range("D14").Value = "4.994"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").Value = "0.6631"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").Value = "30.626.61"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("D17").Value = "0.000007868"
$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").Value = "2.125.91"
$ws.Range("E20").Value = "  -1.64%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").Value = "4.731"
$ws.Range("E22").Value = "  -2.84%  "

$ws.Range("D23").Value = "191.77"
$ws.Range("E23").Value = "  +4.71%  "

$ws.Range("D24").Value = "6.059"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").Value = "9.293"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").Value = "158.65"
$ws.Range("E26").Value = "  +3.40%  "

$ws.Range("D27").Value = "18.28"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").Value = "1.828"
$ws.Range("E28").Value = "  -5.90%  "

$ws.Range("D29").Value = "1.404"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").Value = "4.251"
$ws.Range("E30").Value = "  -2.10%  "

$ws.Range("D31").Value = "0.09030"
$ws.Range("E31").Value = "  +0.69%  "

$ws.Range("D32").Value = "3.931"
$ws.Range("E32").Value = "  -2.87%  "

$ws.Range("D33").Value = "0.05196"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").Value = "0.7305"
$ws.Range("E34").Value = "  -2.36%  "

$ws.Range("E35").Value = "  -5.45%  "

$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("D37").Value = "0.01828"
$ws.Range("E37").Value = "  -5.53%  "

$ws.Range("D38").Value = "2.670"
$ws.Range("E38").Value = "  -2.32%  "

$ws.Range("D39").Value = "0.9226"
$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").Value = "2.050"
$ws.Range("E40").Value = "  -6.05%  "

$ws.Range("D41").Value = "0.4368"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "104.43"
$ws.Range("E42").Value = "  -1.50%  "

$ws.Range("D43").Value = "0.9987"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Value = "5.710"
$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("D45").Value = "0.1334"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").Value = "7.279"
$ws.Range("E46").Value = "  -6.52%  "

$ws.Range("D47").Value = "0.4095"
$ws.Range("E47").Value = "  +5.27%  "

$ws.Range("D48").Value = "0.05828"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").Value = "8.662"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").Value = "1.405"
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("E51").Value = "  -0.22%  "
